# Fruta / hortaliza, semanal
# Update the weekly snapshot of price rows: the Fecha (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M) and
# Precio $/Kg (P) columns are refreshed to the latest data pull for each
# row 2..19 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @{ D = 44894; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    3  = @{ D = 44846; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    4  = @{ D = 44804; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    5  = @{ D = 44810; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
    6  = @{ D = 44839; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    7  = @{ D = 44841; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    8  = @{ D = 44797; J = 60; K = 12000; L = 13000; M = 12500; P = 962 }
    9  = @{ D = 44930; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    10 = @{ D = 44859; J = 30; K = 13000; L = 13000; M = 13000; P = 1000 }
    11 = @{ D = 44895; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    12 = @{ D = 44868; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    13 = @{ D = 44915; J = 50; K = 18000; L = 18000; M = 18000; P = 1385 }
    14 = @{ D = 44943; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    15 = @{ D = 44874; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    16 = @{ D = 44832; J = 60; K = 17000; L = 18000; M = 17500; P = 1346 }
    17 = @{ D = 44922; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    18 = @{ D = 44959; J = 30; K = 19000; L = 19000; M = 19000; P = 1462 }
    19 = @{ D = 44880; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D    # D: Fecha
    $ws.Cells.Item($r, 10).Value = $vals.J    # J: Volumen
    $ws.Cells.Item($r, 11).Value = $vals.K    # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $vals.L    # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $vals.M    # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $vals.P    # P: Precio $/Kg
}
